# Auto-generated Excel COM-interop script
# Applies the diff: inserts a new injury row (Дроздов Иван) into the
# "snapshot" sheet (shifting existing CSKA/Kunlun rows down by one),
# refreshes the K-column scrape timestamps for every existing row, and
# appends the corresponding notification row to the "new_injured" sheet.

$wb = $excel.ActiveWorkbook
$snapshot = $wb.Worksheets.Item("snapshot")
$newInjured = $wb.Worksheets.Item("new_injured")

# --- Refresh scraped_at (column K) timestamps for rows 2-47 (data unchanged) ---
$snapshot.Range("K2").Value = "2025-11-11T07:03:33.328798+00:00"
$snapshot.Range("K3").Value = "2025-11-11T07:03:33.328841+00:00"
$snapshot.Range("K4").Value = "2025-11-11T07:03:33.328868+00:00"
$snapshot.Range("K5").Value = "2025-11-11T07:03:36.128766+00:00"
$snapshot.Range("K6").Value = "2025-11-11T07:03:36.128809+00:00"
$snapshot.Range("K7").Value = "2025-11-11T07:03:36.128830+00:00"
$snapshot.Range("K8").Value = "2025-11-11T07:03:38.959326+00:00"
$snapshot.Range("K9").Value = "2025-11-11T07:03:41.795700+00:00"
$snapshot.Range("K10").Value = "2025-11-11T07:03:41.795734+00:00"
$snapshot.Range("K11").Value = "2025-11-11T07:03:41.795756+00:00"
$snapshot.Range("K12").Value = "2025-11-11T07:03:44.565472+00:00"
$snapshot.Range("K13").Value = "2025-11-11T07:03:44.565507+00:00"
$snapshot.Range("K14").Value = "2025-11-11T07:03:44.565526+00:00"
$snapshot.Range("K15").Value = "2025-11-11T07:03:44.565543+00:00"
$snapshot.Range("K16").Value = "2025-11-11T07:03:49.132440+00:00"
$snapshot.Range("K17").Value = "2025-11-11T07:03:51.852141+00:00"
$snapshot.Range("K18").Value = "2025-11-11T07:03:54.047824+00:00"
$snapshot.Range("K19").Value = "2025-11-11T07:03:54.047863+00:00"
$snapshot.Range("K20").Value = "2025-11-11T07:03:54.047887+00:00"
$snapshot.Range("K21").Value = "2025-11-11T07:03:56.405613+00:00"
$snapshot.Range("K22").Value = "2025-11-11T07:03:59.179203+00:00"
$snapshot.Range("K23").Value = "2025-11-11T07:03:59.179255+00:00"
$snapshot.Range("K24").Value = "2025-11-11T07:04:01.915164+00:00"
$snapshot.Range("K25").Value = "2025-11-11T07:04:01.915198+00:00"
$snapshot.Range("K26").Value = "2025-11-11T07:04:01.915219+00:00"
$snapshot.Range("K27").Value = "2025-11-11T07:04:04.148967+00:00"
$snapshot.Range("K28").Value = "2025-11-11T07:04:04.148997+00:00"
$snapshot.Range("K29").Value = "2025-11-11T07:04:04.149017+00:00"
$snapshot.Range("K30").Value = "2025-11-11T07:04:04.149035+00:00"
$snapshot.Range("K31").Value = "2025-11-11T07:04:04.149052+00:00"
$snapshot.Range("K32").Value = "2025-11-11T07:04:06.946952+00:00"
$snapshot.Range("K33").Value = "2025-11-11T07:04:06.946984+00:00"
$snapshot.Range("K34").Value = "2025-11-11T07:04:09.641582+00:00"
$snapshot.Range("K35").Value = "2025-11-11T07:04:09.641616+00:00"
$snapshot.Range("K36").Value = "2025-11-11T07:04:09.641638+00:00"
$snapshot.Range("K37").Value = "2025-11-11T07:04:11.974401+00:00"
$snapshot.Range("K38").Value = "2025-11-11T07:04:11.974439+00:00"
$snapshot.Range("K39").Value = "2025-11-11T07:04:11.974460+00:00"
$snapshot.Range("K40").Value = "2025-11-11T07:04:14.277434+00:00"
$snapshot.Range("K41").Value = "2025-11-11T07:04:14.277468+00:00"
$snapshot.Range("K42").Value = "2025-11-11T07:04:14.277539+00:00"
$snapshot.Range("K43").Value = "2025-11-11T07:04:14.277568+00:00"
$snapshot.Range("K44").Value = "2025-11-11T07:04:14.277588+00:00"
$snapshot.Range("K45").Value = "2025-11-11T07:04:14.277605+00:00"
$snapshot.Range("K46").Value = "2025-11-11T07:04:16.556140+00:00"
$snapshot.Range("K47").Value = "2025-11-11T07:04:16.556171+00:00"

# --- Shift CSKA/Kunlun block down by one row and insert the new injury ---
# (rows 48-52 -> 49-53, unchanged values, refreshed K timestamps; new row 48 inserted)
# Row 48
$snapshot.Range("A48").Value = "ЦСК"
$snapshot.Range("B48").Value = "ЦСКА"
$snapshot.Range("C48").Value = "cska"
$snapshot.Range("D48").Value = "Дроздов Иван"
$snapshot.Range("E48").Value = "19"
$snapshot.Range("F48").Value = "нападающий"
$snapshot.Range("G48").Value = "30752"
$snapshot.Range("H48").Value = "1369_ЦСК_дроздовиван"
$snapshot.Range("I48").Value = "injured_active"
$snapshot.Range("J48").Value = "https://www.khl.ru/clubs/cska/team/"
$snapshot.Range("K48").Value = "2025-11-11T07:04:21.279309+00:00"
# Row 49
$snapshot.Range("A49").Value = "ЦСК"
$snapshot.Range("B49").Value = "ЦСКА"
$snapshot.Range("C49").Value = "cska"
$snapshot.Range("D49").Value = "Моисеев Данила"
$snapshot.Range("E49").Value = "93"
$snapshot.Range("F49").Value = "нападающий"
$snapshot.Range("G49").Value = "23931"
$snapshot.Range("H49").Value = "1369_ЦСК_моисеевданила"
$snapshot.Range("I49").Value = "injured_active"
$snapshot.Range("J49").Value = "https://www.khl.ru/clubs/cska/team/"
$snapshot.Range("K49").Value = "2025-11-11T07:04:21.279336+00:00"
# Row 50
$snapshot.Range("A50").Value = "ЦСК"
$snapshot.Range("B50").Value = "ЦСКА"
$snapshot.Range("C50").Value = "cska"
$snapshot.Range("D50").Value = "Саморуков Дмитрий"
$snapshot.Range("E50").Value = "5"
$snapshot.Range("F50").Value = "защитник"
$snapshot.Range("G50").Value = "24005"
$snapshot.Range("H50").Value = "1369_ЦСК_саморуковдмитрий"
$snapshot.Range("I50").Value = "injured_active"
$snapshot.Range("J50").Value = "https://www.khl.ru/clubs/cska/team/"
$snapshot.Range("K50").Value = "2025-11-11T07:04:21.279354+00:00"
# Row 51
$snapshot.Range("A51").Value = "ЦСК"
$snapshot.Range("B51").Value = "ЦСКА"
$snapshot.Range("C51").Value = "cska"
$snapshot.Range("D51").Value = "Уильямс Колби"
$snapshot.Range("E51").Value = "22"
$snapshot.Range("F51").Value = "защитник"
$snapshot.Range("G51").Value = "41896"
$snapshot.Range("H51").Value = "1369_ЦСК_уильямсколби"
$snapshot.Range("I51").Value = "injured_active"
$snapshot.Range("J51").Value = "https://www.khl.ru/clubs/cska/team/"
$snapshot.Range("K51").Value = "2025-11-11T07:04:21.279373+00:00"
# Row 52
$snapshot.Range("A52").Value = "ШДР"
$snapshot.Range("B52").Value = "Драконы"
$snapshot.Range("C52").Value = "kunlun"
$snapshot.Range("D52").Value = "Гроло Жереми"
$snapshot.Range("E52").Value = "75"
$snapshot.Range("F52").Value = "защитник"
$snapshot.Range("G52").Value = "45343"
$snapshot.Range("H52").Value = "1369_ШДР_гроложереми"
$snapshot.Range("I52").Value = "injured_active"
$snapshot.Range("J52").Value = "https://www.khl.ru/clubs/kunlun/team/"
$snapshot.Range("K52").Value = "2025-11-11T07:04:23.979096+00:00"
# Row 53
$snapshot.Range("A53").Value = "ШДР"
$snapshot.Range("B53").Value = "Драконы"
$snapshot.Range("C53").Value = "kunlun"
$snapshot.Range("D53").Value = "Саттер Райли"
$snapshot.Range("E53").Value = "14"
$snapshot.Range("F53").Value = "нападающий"
$snapshot.Range("G53").Value = "45491"
$snapshot.Range("H53").Value = "1369_ШДР_саттеррайли"
$snapshot.Range("I53").Value = "injured_active"
$snapshot.Range("J53").Value = "https://www.khl.ru/clubs/kunlun/team/"
$snapshot.Range("K53").Value = "2025-11-11T07:04:23.979125+00:00"

# --- Append the new_injured notification row ---
# Row 2
$newInjured.Range("A2").Value = "ЦСК"
$newInjured.Range("B2").Value = "ЦСКА"
$newInjured.Range("C2").Value = "Дроздов Иван"
$newInjured.Range("D2").Value = "1369_ЦСК_дроздовиван"
$newInjured.Range("E2").Value = "INJURED_NEW"
$newInjured.Range("F2").Value = "2025-11-11T15:04:24.488351+08:00"
$newInjured.Range("G2").Value = "2025-11-11"
